$wb = $excel.ActiveWorkbook

# --- Sheet 1: journalVoucherDetails -----------------------------------
$ws1 = $wb.Worksheets.Item("journalVoucherDetails")

# New column H header: accountCode3
$ws1.Cells.Item(1, 8).Value = "accountCode3"

# New row 6: the "remittance" data row (accountCode3 = 3502002)
$ws1.Cells.Item(6, 1).Value = "remittance"
$ws1.Cells.Item(6, 2).Value = "03/01/2016"
$ws1.Cells.Item(6, 3).Value = "Expense"
$ws1.Cells.Item(6, 4).Value = "2101001"
$ws1.Cells.Item(6, 5).Value = "3501003"
$ws1.Cells.Item(6, 6).Value = "PUBLIC HEALTH AND SANITATION"
$ws1.Cells.Item(6, 7).Value = "Public Health"
$ws1.Cells.Item(6, 8).Value = 3502002

# --- Sheet 2: financialBankDetails -------------------------------------
$ws2 = $wb.Worksheets.Item("financialBankDetails")

# New row 3: second bank/account entry
$ws2.Cells.Item(3, 1).Value = "SBI1"
$ws2.Cells.Item(3, 2).Value = "KOTAK MAHINDRA BANK Ucon Plaza Kurnool"
$ws2.Cells.Item(3, 3).Value = "4502205--311010192115--KOTAK MAHINDRA BANK"

# --- View state: activate financialBankDetails, update selections -----
$ws1.Range("H6").Select()
$ws2.Activate()
$ws2.Range("C3").Select()
